$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix typo: "inspire" -> "inspired"
$ws.Range("D3").Value = 'As a User, I want to get an overview on the Mainpage on how this S.T.A.L.K.E.R inspired server looks like'

# Update selection/view to reflect where the edit was made
$ws.Activate()
$ws.Range("D3").Select()
